$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("Epoch Accuracy") value corrections -------------------------
# Re-run of the eval cell (new random seed / weight init -> new accuracy
# numbers for epochs 3..116 plus the final overall accuracy in B118).
$accuracyUpdates = [ordered]@{
    5  = 0.3125
    7  = 0.3125
    8  = 0.296875
    9  = 0.296875
    10  = 0.296875
    11  = 0.28125
    12  = 0.390625
    13  = 0.390625
    14  = 0.328125
    15  = 0.265625
    16  = 0.328125
    17  = 0.34375
    18  = 0.40625
    19  = 0.34375
    20  = 0.328125
    21  = 0.3125
    22  = 0.3125
    23  = 0.28125
    24  = 0.28125
    25  = 0.296875
    26  = 0.296875
    27  = 0.296875
    28  = 0.296875
    29  = 0.28125
    30  = 0.28125
    31  = 0.28125
    32  = 0.28125
    33  = 0.28125
    34  = 0.28125
    35  = 0.296875
    36  = 0.296875
    37  = 0.296875
    38  = 0.296875
    39  = 0.296875
    40  = 0.296875
    41  = 0.296875
    42  = 0.296875
    43  = 0.296875
    44  = 0.296875
    45  = 0.296875
    46  = 0.296875
    47  = 0.296875
    48  = 0.296875
    49  = 0.296875
    50  = 0.296875
    51  = 0.296875
    52  = 0.296875
    53  = 0.296875
    54  = 0.296875
    55  = 0.296875
    56  = 0.296875
    57  = 0.296875
    58  = 0.296875
    59  = 0.296875
    60  = 0.296875
    61  = 0.296875
    62  = 0.296875
    63  = 0.296875
    64  = 0.296875
    65  = 0.296875
    66  = 0.296875
    67  = 0.296875
    68  = 0.296875
    69  = 0.296875
    70  = 0.296875
    71  = 0.296875
    72  = 0.296875
    73  = 0.296875
    74  = 0.296875
    75  = 0.296875
    76  = 0.296875
    77  = 0.296875
    78  = 0.296875
    79  = 0.296875
    80  = 0.296875
    81  = 0.296875
    82  = 0.296875
    83  = 0.296875
    84  = 0.296875
    85  = 0.296875
    86  = 0.296875
    87  = 0.296875
    88  = 0.296875
    89  = 0.296875
    90  = 0.296875
    91  = 0.296875
    92  = 0.296875
    93  = 0.296875
    94  = 0.296875
    95  = 0.296875
    96  = 0.296875
    97  = 0.296875
    98  = 0.296875
    99  = 0.296875
    100  = 0.296875
    101  = 0.296875
    102  = 0.296875
    103  = 0.3125
    104  = 0.296875
    105  = 0.3125
    106  = 0.234375
    107  = 0.15625
    108  = 0.296875
    109  = 0.203125
    110  = 0.203125
    111  = 0.3125
    113  = 0.296875
    114  = 0.1875
    115  = 0.140625
    117  = 0.265625
    118  = 0.2131147540983606
}
foreach ($row in $accuracyUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $accuracyUpdates[$row]
}

# --- Column A inline-string repr refresh (rows 102-118) --------------------
# The object repr embeds the python process's memory address; it changed
# because the notebook kernel was restarted for this run.
$newRepr = "<__main__.DisplayOutputs object at 0x7f7d10189910>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}

# --- Sheet view ---------------------------------------------------------
# Drop the frozen scroll position (topLeftCell) and leave the cursor on E11
# while the whole sheet remains selected.
$ws.Range("E11").Select()
$ws.Range("A1:XFD1048576").Select()
